# Update the "Förändrad" (Changed) date column (C) for rows 2-8
# from 2023-10-13 (serial 45212) to 2023-10-22 (serial 45221).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    $current = $cell.Value()
    $serial = $current.ToOADate()
    if ($serial -eq 45212) {
        $cell.Value = 45221
    }
}
